$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (Dom Perignon)
# ---------------------------------------------------------------------------
$ws.Range("AC2").Value = "750 ML"
$ws.Range("AU2").Value = "2025-03-29T02:10:07.293Z"
$ws.Range("AV2").ClearContents()

# ---------------------------------------------------------------------------
# Row 5 (Grey Goose Vodka)
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "Grey Goose Vodka 34"

# ---------------------------------------------------------------------------
# Row 24 (Late-Night Jazz / SKU-606-998)
# ---------------------------------------------------------------------------
$ws.Range("AB24").Value = "34ml"
$ws.Range("AC24").ClearContents()
$ws.Range("AG24").Value = 34
$ws.Range("AU24").Value = "2025-03-28T15:15:57.426Z"

# ---------------------------------------------------------------------------
# Row 25 - fully replaced with "test" placeholder data (id 607 kept)
# ---------------------------------------------------------------------------
$ws.Rows.Item(25).ClearContents()

$ws.Range("A25").Value = 607
$ws.Range("B25").Value = "SKU-606-999"
$ws.Range("C25").Value = "test"
$ws.Range("D25").Value = "test"
$ws.Range("E25").Value = "test"
$ws.Range("F25").Value = "test category 34"
$ws.Range("G25").Value = "test"
$ws.Range("H25").Value = "test"
$ws.Range("I25").Value = "test"
$ws.Range("J25").Value = "test"
$ws.Range("K25").Value = "test"
$ws.Range("L25").Value = "test"
$ws.Range("M25").Value = "test"
$ws.Range("N25").Value = "test"
$ws.Range("O25").Value = "test"
$ws.Range("P25").Value = "test"
$ws.Range("Q25").Value = 56
$ws.Range("R25").Value = 56
$ws.Range("S25").Value = 56
$ws.Range("T25").Value = "https://res.cloudinary.com/dc3hqcovg/image/upload/v1743174957/gkmedfirj67byqafpo36.png"
$ws.Range("AB25").Value = "76ml"
$ws.Range("AG25").Value = 76
$ws.Range("AV25").Value = "test category"

# ---------------------------------------------------------------------------
# Row 26 (Five Stage Pipeline -> test product)
# ---------------------------------------------------------------------------
$ws.Range("B26").Value = "SKU-608-717"
$ws.Range("C26").Value = "test product"
$ws.Range("D26").Value = "test product"
$ws.Range("F26").Value = "beer"
$ws.Range("G26").Value = "whiskey"
$ws.Range("T26").Value = "https://res.cloudinary.com/dc3hqcovg/image/upload/v1743226484/vzcjhlpqnvz4i0x2uw1j.svg"
$ws.Range("AG26").Value = 34
$ws.Range("AT26").Value = "https://res.cloudinary.com/dc3hqcovg/image/upload/v1743226484/vzcjhlpqnvz4i0x2uw1j.svg"
$ws.Range("AU26").Value = "2025-03-29T05:34:46.251Z"
$ws.Range("AV26").Value = "beer"
$ws.Range("AX26").Value = "2025-03-29T05:34:46.251Z"
